# Update "Datos actualizados" timestamp cell (A1) on the "Ciudades" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 20:16"

# Row 15 corresponds to "Asturias": Recuperados (D15) 329 -> 328, Muertes (E15) 3 -> 4
$ws.Range("D15").Value = 328
$ws.Range("E15").Value = 4

# Row 22 corresponds to "Murcia": Recuperados (D22) 214 -> 213, Muertes (E22) 0 -> 1
$ws.Range("D22").Value = 213
$ws.Range("E22").Value = 1
